$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21 (shifts existing rows 21..112 down to 22..113)
$ws.Rows("21").Insert()

# Populate the new row 21 with the new daily record
$ws.Range("A21").Value = 5
$ws.Range("B21").Value = "Macroferia Regional de Talca"
$ws.Range("C21").Value = "Maule"
$ws.Range("D21").Value = 44613
$ws.Range("E21").Value = 7
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100108
$ws.Range("H21").Value = "Tropicales y subtropicales"
$ws.Range("I21").Value = 100108002
$ws.Range("J21").Value = "Mango"
$ws.Range("K21").Value = "Sin especificar"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 200
$ws.Range("N21").Value = 7000
$ws.Range("O21").Value = 7000
$ws.Range("P21").Value = 7000
$ws.Range("Q21").Value = "$/bandeja 4 kilos"
$ws.Range("R21").Value = "Perú"
$ws.Range("S21").Value = 1750
$ws.Range("T21").Value = 4
